$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 11:43"

# Update country names (column A) where rows were re-sorted/re-ordered
$ws.Cells.Item(23, 1).Value = "Indonesia"
$ws.Cells.Item(24, 1).Value = "Turquia"
$ws.Cells.Item(36, 1).Value = "Polonia"
$ws.Cells.Item(37, 1).Value = "Catar"
$ws.Cells.Item(44, 1).Value = "Oman"
$ws.Cells.Item(45, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(91, 1).Value = "Eslovaquia"
$ws.Cells.Item(92, 1).Value = "Costa de Marfil"
$ws.Cells.Item(150, 1).Value = "Letonia"
$ws.Cells.Item(151, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(152, 1).Value = "Principado de Andorra"

# Update numeric data columns (B..H) per row
$ws.Cells.Item(23, 2).Value = 336716
$ws.Cells.Item(23, 3).Value = 3267
$ws.Cells.Item(23, 4).Value = 258519
$ws.Cells.Item(23, 5).Value = 66262
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 91
$ws.Cells.Item(23, 8).Value = 11935

$ws.Cells.Item(24, 2).Value = 335533
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 294357
$ws.Cells.Item(24, 5).Value = 32339
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 8837

$ws.Cells.Item(25, 2).Value = 326482
$ws.Cells.Item(25, 3).Value = 191
$ws.Cells.Item(25, 4).Value = 274700
$ws.Cells.Item(25, 5).Value = 42080
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 9702

$ws.Cells.Item(36, 2).Value = 130210
$ws.Cells.Item(36, 3).Value = 4394
$ws.Cells.Item(36, 4).Value = 81201
$ws.Cells.Item(36, 5).Value = 45970
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 35
$ws.Cells.Item(36, 8).Value = 3039

$ws.Cells.Item(37, 2).Value = 127985
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 124978
$ws.Cells.Item(37, 5).Value = 2787
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 220

$ws.Cells.Item(40, 4).Value = 56203
$ws.Cells.Item(40, 5).Value = 59920

$ws.Cells.Item(44, 2).Value = 106575
$ws.Cells.Item(44, 3).Value = 685
$ws.Cells.Item(44, 4).Value = 93222
$ws.Cells.Item(44, 5).Value = 12307
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 8
$ws.Cells.Item(44, 8).Value = 1046

$ws.Cells.Item(45, 2).Value = 106229
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 97284
$ws.Cells.Item(45, 5).Value = 8500
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 445

$ws.Cells.Item(57, 5).Value = 4088
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 277

$ws.Cells.Item(64, 2).Value = 56298
$ws.Cells.Item(64, 3).Value = 979
$ws.Cells.Item(64, 4).Value = 44065
$ws.Cells.Item(64, 5).Value = 11378
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 4
$ws.Cells.Item(64, 8).Value = 855

$ws.Cells.Item(91, 2).Value = 20355
$ws.Cells.Item(91, 3).Value = 504
$ws.Cells.Item(91, 4).Value = 6031
$ws.Cells.Item(91, 5).Value = 14263
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 61

$ws.Cells.Item(92, 2).Value = 20154
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 19798
$ws.Cells.Item(92, 5).Value = 236
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 120

$ws.Cells.Item(101, 5).Value = 5642
$ws.Cells.Item(101, 7).Value = 7
$ws.Cells.Item(101, 8).Value = 92

$ws.Cells.Item(121, 2).Value = 6248
$ws.Cells.Item(121, 3).Value = 126
$ws.Cells.Item(121, 4).Value = 2793
$ws.Cells.Item(121, 5).Value = 3352

$ws.Cells.Item(127, 2).Value = 5194
$ws.Cells.Item(127, 3).Value = 11
$ws.Cells.Item(127, 4).Value = 4921
$ws.Cells.Item(127, 5).Value = 168

$ws.Cells.Item(135, 2).Value = 4791
$ws.Cells.Item(135, 3).Value = 39
$ws.Cells.Item(135, 4).Value = 3317
$ws.Cells.Item(135, 5).Value = 1461

$ws.Cells.Item(150, 2).Value = 2765
$ws.Cells.Item(150, 3).Value = 95
$ws.Cells.Item(150, 4).Value = 1325
$ws.Cells.Item(150, 5).Value = 1400
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 40

$ws.Cells.Item(151, 2).Value = 2754
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 2019
$ws.Cells.Item(151, 5).Value = 725
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 10

$ws.Cells.Item(152, 2).Value = 2696
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 1814
$ws.Cells.Item(152, 5).Value = 827
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 55

$ws.Cells.Item(158, 2).Value = 2280
$ws.Cells.Item(158, 3).Value = 9
$ws.Cells.Item(158, 4).Value = 1554
$ws.Cells.Item(158, 5).Value = 663
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 2
$ws.Cells.Item(158, 8).Value = 63

